$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column is BM = column 65 (one past the previous last column BL = 64)
$newCol = 65

# Match the new column's width to the preceding column (BL), which carries
# width="12" customWidth="1" in the OOXML <cols> block.
$ws.Columns.Item($newCol).ColumnWidth = $ws.Columns.Item($newCol - 1).ColumnWidth

# --- Header cell BM1: "2024/11/12" as literal text (same look as B1:BL1) ---
# A leading apostrophe forces the date-look-alike string to be kept as text
# instead of Excel auto-converting it to a date serial number.
$headerCell = $ws.Cells.Item(1, $newCol)
$headerCell.Value = "'2024/11/12"
# Re-apply the same cell format (font, no fill, General number format) as the
# neighbouring header cell so BM1 matches B1:BL1 exactly.
$ws.Range("B1").Copy() | Out-Null
$headerCell.PasteSpecial(-4122) | Out-Null

# --- Data rows 2-53: composite-probability values for 2024/11/12 ---
# Style 1 = default (font only, no fill)      -> value >= ~140
# Style 2 = yellow fill (font)                -> value <  ~125
# Style 3 = light-blue fill (font)             -> value in between
$styleTemplateCells = @{
    1 = $ws.Range("B2")
    2 = $ws.Range("D2")
    3 = $ws.Range("N2")
}

$newColumnData = @(
    @{ Row = 2; Style = 2; Value = 124.1 },
    @{ Row = 3; Style = 1; Value = 172.4 },
    @{ Row = 4; Style = 1; Value = 242.8 },
    @{ Row = 5; Style = 1; Value = 156.3 },
    @{ Row = 6; Style = 2; Value = 123.9 },
    @{ Row = 7; Style = 1; Value = 627 },
    @{ Row = 8; Style = 1; Value = 163.7 },
    @{ Row = 9; Style = 1; Value = 213.6 },
    @{ Row = 10; Style = 1; Value = 186.2 },
    @{ Row = 11; Style = 1; Value = 303.6 },
    @{ Row = 12; Style = 1; Value = 158.7 },
    @{ Row = 13; Style = 1; Value = 140 },
    @{ Row = 14; Style = 2; Value = 107.4 },
    @{ Row = 15; Style = 1; Value = 146.8 },
    @{ Row = 16; Style = 2; Value = 117.8 },
    @{ Row = 17; Style = 2; Value = 111.7 },
    @{ Row = 18; Style = 1; Value = 151.4 },
    @{ Row = 19; Style = 1; Value = 165.2 },
    @{ Row = 20; Style = 3; Value = 136.9 },
    @{ Row = 21; Style = 2; Value = 121.3 },
    @{ Row = 22; Style = 1; Value = 163.4 },
    @{ Row = 23; Style = 1; Value = 165.7 },
    @{ Row = 24; Style = 1; Value = 168.1 },
    @{ Row = 25; Style = 3; Value = 130 },
    @{ Row = 26; Style = 3; Value = 131.8 },
    @{ Row = 27; Style = 1; Value = 143.6 },
    @{ Row = 28; Style = 1; Value = 140.1 },
    @{ Row = 29; Style = 1; Value = 179.8 },
    @{ Row = 30; Style = 1; Value = 207.3 },
    @{ Row = 31; Style = 1; Value = 180.4 },
    @{ Row = 32; Style = 3; Value = 137 },
    @{ Row = 33; Style = 1; Value = 144.9 },
    @{ Row = 34; Style = 1; Value = 151.2 },
    @{ Row = 35; Style = 3; Value = 126.3 },
    @{ Row = 36; Style = 1; Value = 183.2 },
    @{ Row = 37; Style = 3; Value = 138.7 },
    @{ Row = 38; Style = 1; Value = 147 },
    @{ Row = 39; Style = 1; Value = 155.4 },
    @{ Row = 40; Style = 2; Value = 122.1 },
    @{ Row = 41; Style = 3; Value = 126.2 },
    @{ Row = 42; Style = 1; Value = 274 },
    @{ Row = 43; Style = 1; Value = 181.1 },
    @{ Row = 44; Style = 2; Value = 119.7 },
    @{ Row = 45; Style = 1; Value = 182 },
    @{ Row = 46; Style = 1; Value = 172.3 },
    @{ Row = 47; Style = 1; Value = 147.7 },
    @{ Row = 48; Style = 1; Value = 153.7 },
    @{ Row = 49; Style = 1; Value = 146.6 },
    @{ Row = 50; Style = 3; Value = 133.6 },
    @{ Row = 51; Style = 1; Value = 238.4 },
    @{ Row = 52; Style = 1; Value = 211.4 },
    @{ Row = 53; Style = 1; Value = 224.3 }

)

foreach ($item in $newColumnData) {
    $cell = $ws.Cells.Item($item.Row, $newCol)
    $cell.Value = $item.Value
    $styleTemplateCells[$item.Style].Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false
